{"js": "// Update the embedded Java exception stack-trace text that appears twice in\n// this bookmark-runtime-exception fixture (M2Doc 3.1.1 -> 3.2.0 regeneration):\n// line numbers shifted in M2DocEvaluator.java / M2DocUtils.java /\n// AbstractTemplatesTestSuite.java, the generated-accessor class name changed,\n// and the tail of the stack (Eclipse JDT JUnit launcher frames) was replaced\n// by the Maven Surefire / Tycho / Equinox launcher frames.\nconst replacements = [\n  [\"M2DocEvaluator.java:1607)\", \"M2DocEvaluator.java:1696)\", 1],\n  [\"M2DocEvaluator.java:1242)\", \"M2DocEvaluator.java:1331)\", 3],\n  [\"M2DocEvaluator.java:1467)\", \"M2DocEvaluator.java:1556)\", 1],\n  [\"M2DocEvaluator.java:297)\", \"M2DocEvaluator.java:301)\", 1],\n  [\"M2DocEvaluator.java:282)\", \"M2DocEvaluator.java:286)\", 1],\n  [\"M2DocUtils.java:845)\", \"M2DocUtils.java:853)\", 1],\n  [\"AbstractTemplatesTestSuite.java:514)\", \"AbstractTemplatesTestSuite.java:518)\", 1],\n  [\"AbstractTemplatesTestSuite.java:421)\", \"AbstractTemplatesTestSuite.java:414)\", 1],\n  [\"GeneratedMethodAccessor73\", \"GeneratedMethodAccessor5\", 1],\n  [\"\\tat org.eclipse.jdt.internal.junit4.runner.JUnit4TestReference.run(JUnit4TestReference.java:86)\\n\\tat org.eclipse.jdt.internal.junit.runner.TestExecution.run(TestExecution.java:38)\\n\\tat org.eclipse.jdt.internal.junit.runner.RemoteTestRunner.runTests(RemoteTestRunner.java:538)\\n\\tat org.eclipse.jdt.internal.junit.runner.RemoteTestRunner.runTests(RemoteTestRunner.java:760)\\n\\tat org.eclipse.jdt.internal.junit.runner.RemoteTestRunner.run(RemoteTestRunner.java:460)\\n\\tat org.eclipse.jdt.internal.junit.runner.RemoteTestRunner.main(RemoteTestRunner.java:206)\", \"\\tat org.apache.maven.surefire.junit4.JUnit4Provider.execute(JUnit4Provider.java:365)\\n\\tat org.apache.maven.surefire.junit4.JUnit4Provider.executeWithRerun(JUnit4Provider.java:273)\\n\\tat org.apache.maven.surefire.junit4.JUnit4Provider.executeTestSet(JUnit4Provider.java:238)\\n\\tat org.apache.maven.surefire.junit4.JUnit4Provider.invoke(JUnit4Provider.java:159)\\n\\tat sun.reflect.NativeMethodAccessorImpl.invoke0(Native Method)\\n\\tat sun.reflect.NativeMethodAccessorImpl.invoke(NativeMethodAccessorImpl.java:62)\\n\\tat sun.reflect.DelegatingMethodAccessorImpl.invoke(DelegatingMethodAccessorImpl.java:43)\\n\\tat java.lang.reflect.Method.invoke(Method.java:498)\\n\\tat org.apache.maven.surefire.util.ReflectionUtils.invokeMethodWithArray2(ReflectionUtils.java:206)\\n\\tat org.apache.maven.surefire.booter.ProviderFactory$ProviderProxy.invoke(ProviderFactory.java:161)\\n\\tat org.apache.maven.surefire.booter.ProviderFactory.invokeProvider(ProviderFactory.java:84)\\n\\tat org.eclipse.tycho.surefire.osgibooter.OsgiSurefireBooter.run(OsgiSurefireBooter.java:113)\\n\\tat org.eclipse.tycho.surefire.osgibooter.HeadlessTestApplication.run(HeadlessTestApplication.java:21)\\n\\tat sun.reflect.NativeMethodAccessorImpl.invoke0(Native Method)\\n\\tat sun.reflect.NativeMethodAccessorImpl.invoke(NativeMethodAccessorImpl.java:62)\\n\\tat sun.reflect.DelegatingMethodAccessorImpl.invoke(DelegatingMethodAccessorImpl.java:43)\\n\\tat java.lang.reflect.Method.invoke(Method.java:498)\\n\\tat org.eclipse.equinox.internal.app.EclipseAppContainer.callMethodWithException(EclipseAppContainer.java:593)\\n\\tat org.eclipse.equinox.internal.app.EclipseAppHandle.run(EclipseAppHandle.java:205)\\n\\tat org.eclipse.core.runtime.internal.adaptor.EclipseAppLauncher.runApplication(EclipseAppLauncher.java:137)\\n\\tat org.eclipse.core.runtime.internal.adaptor.EclipseAppLauncher.start(EclipseAppLauncher.java:107)\\n\\tat org.eclipse.core.runtime.adaptor.EclipseStarter.run(EclipseStarter.java:401)\\n\\tat org.eclipse.core.runtime.adaptor.EclipseStarter.run(EclipseStarter.java:255)\\n\\tat sun.reflect.NativeMethodAccessorImpl.invoke0(Native Method)\\n\\tat sun.reflect.NativeMethodAccessorImpl.invoke(NativeMethodAccessorImpl.java:62)\\n\\tat sun.reflect.DelegatingMethodAccessorImpl.invoke(DelegatingMethodAccessorImpl.java:43)\\n\\tat java.lang.reflect.Method.invoke(Method.java:498)\\n\\tat org.eclipse.equinox.launcher.Main.invokeFramework(Main.java:657)\\n\\tat org.eclipse.equinox.launcher.Main.basicRun(Main.java:594)\\n\\tat org.eclipse.equinox.launcher.Main.run(Main.java:1447)\\n\\tat org.eclipse.equinox.launcher.Main.main(Main.java:1420)\", 1],\n];\n\nconst body = context.document.body;\n\nfor (const [search, replace, expectedCount] of replacements) {\n  const results = body.search(search, { matchCase: true, matchWildcards: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length !== expectedCount) {\n    throw new Error(\n      `Expected ${expectedCount} match(es) for ${JSON.stringify(search.substring(0, 60))}, found ${results.items.length}`\n    );\n  }\n\n  for (const item of results.items) {\n    item.insertText(replace, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Update the embedded Java exception stack-trace text that appears twice in\n# this bookmark-runtime-exception fixture (M2Doc 3.1.1 -> 3.2.0 regeneration):\n# line numbers shifted in M2DocEvaluator.java / M2DocUtils.java /\n# AbstractTemplatesTestSuite.java, the generated-accessor class name changed,\n# and the tail of the stack (Eclipse JDT JUnit launcher frames) was replaced\n# by the Maven Surefire / Tycho / Equinox launcher frames.\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n  ,@(\"M2DocEvaluator.java:1607)\", \"M2DocEvaluator.java:1696)\", 1)\n  ,@(\"M2DocEvaluator.java:1242)\", \"M2DocEvaluator.java:1331)\", 3)\n  ,@(\"M2DocEvaluator.java:1467)\", \"M2DocEvaluator.java:1556)\", 1)\n  ,@(\"M2DocEvaluator.java:297)\", \"M2DocEvaluator.java:301)\", 1)\n  ,@(\"M2DocEvaluator.java:282)\", \"M2DocEvaluator.java:286)\", 1)\n  ,@(\"M2DocUtils.java:845)\", \"M2DocUtils.java:853)\", 1)\n  ,@(\"AbstractTemplatesTestSuite.java:514)\", \"AbstractTemplatesTestSuite.java:518)\", 1)\n  ,@(\"AbstractTemplatesTestSuite.java:421)\", \"AbstractTemplatesTestSuite.java:414)\", 1)\n  ,@(\"GeneratedMethodAccessor73\", \"GeneratedMethodAccessor5\", 1)\n  ,@(\"`tat org.eclipse.jdt.internal.junit4.runner.JUnit4TestReference.run(JUnit4TestReference.java:86)`n`tat org.eclipse.jdt.internal.junit.runner.TestExecution.run(TestExecution.java:38)`n`tat org.eclipse.jdt.internal.junit.runner.RemoteTestRunner.runTests(RemoteTestRunner.java:538)`n`tat org.eclipse.jdt.internal.junit.runner.RemoteTestRunner.runTests(RemoteTestRunner.java:760)`n`tat org.eclipse.jdt.internal.junit.runner.RemoteTestRunner.run(RemoteTestRunner.java:460)`n`tat org.eclipse.jdt.internal.junit.runner.RemoteTestRunner.main(RemoteTestRunner.java:206)\", \"`tat org.apache.maven.surefire.junit4.JUnit4Provider.execute(JUnit4Provider.java:365)`n`tat org.apache.maven.surefire.junit4.JUnit4Provider.executeWithRerun(JUnit4Provider.java:273)`n`tat org.apache.maven.surefire.junit4.JUnit4Provider.executeTestSet(JUnit4Provider.java:238)`n`tat org.apache.maven.surefire.junit4.JUnit4Provider.invoke(JUnit4Provider.java:159)`n`tat sun.reflect.NativeMethodAccessorImpl.invoke0(Native Method)`n`tat sun.reflect.NativeMethodAccessorImpl.invoke(NativeMethodAccessorImpl.java:62)`n`tat sun.reflect.DelegatingMethodAccessorImpl.invoke(DelegatingMethodAccessorImpl.java:43)`n`tat java.lang.reflect.Method.invoke(Method.java:498)`n`tat org.apache.maven.surefire.util.ReflectionUtils.invokeMethodWithArray2(ReflectionUtils.java:206)`n`tat org.apache.maven.surefire.booter.ProviderFactory`$ProviderProxy.invoke(ProviderFactory.java:161)`n`tat org.apache.maven.surefire.booter.ProviderFactory.invokeProvider(ProviderFactory.java:84)`n`tat org.eclipse.tycho.surefire.osgibooter.OsgiSurefireBooter.run(OsgiSurefireBooter.java:113)`n`tat org.eclipse.tycho.surefire.osgibooter.HeadlessTestApplication.run(HeadlessTestApplication.java:21)`n`tat sun.reflect.NativeMethodAccessorImpl.invoke0(Native Method)`n`tat sun.reflect.NativeMethodAccessorImpl.invoke(NativeMethodAccessorImpl.java:62)`n`tat sun.reflect.DelegatingMethodAccessorImpl.invoke(DelegatingMethodAccessorImpl.java:43)`n`tat java.lang.reflect.Method.invoke(Method.java:498)`n`tat org.eclipse.equinox.internal.app.EclipseAppContainer.callMethodWithException(EclipseAppContainer.java:593)`n`tat org.eclipse.equinox.internal.app.EclipseAppHandle.run(EclipseAppHandle.java:205)`n`tat org.eclipse.core.runtime.internal.adaptor.EclipseAppLauncher.runApplication(EclipseAppLauncher.java:137)`n`tat org.eclipse.core.runtime.internal.adaptor.EclipseAppLauncher.start(EclipseAppLauncher.java:107)`n`tat org.eclipse.core.runtime.adaptor.EclipseStarter.run(EclipseStarter.java:401)`n`tat org.eclipse.core.runtime.adaptor.EclipseStarter.run(EclipseStarter.java:255)`n`tat sun.reflect.NativeMethodAccessorImpl.invoke0(Native Method)`n`tat sun.reflect.NativeMethodAccessorImpl.invoke(NativeMethodAccessorImpl.java:62)`n`tat sun.reflect.DelegatingMethodAccessorImpl.invoke(DelegatingMethodAccessorImpl.java:43)`n`tat java.lang.reflect.Method.invoke(Method.java:498)`n`tat org.eclipse.equinox.launcher.Main.invokeFramework(Main.java:657)`n`tat org.eclipse.equinox.launcher.Main.basicRun(Main.java:594)`n`tat org.eclipse.equinox.launcher.Main.run(Main.java:1447)`n`tat org.eclipse.equinox.launcher.Main.main(Main.java:1420)\", 1)\n)\n\nforeach ($pair in $replacements) {\n  $search = $pair[0]\n  $replace = $pair[1]\n  $expectedCount = $pair[2]\n\n  $count = 0\n  $find = $d.Content.Find\n  $find.Text = $search\n  while ($find.Execute()) {\n    $count = $count + 1\n  }\n\n  if ($count -ne $expectedCount) {\n    throw \"Expected $expectedCount match(es) for '$($search.Substring(0, [Math]::Min(60, $search.Length))) ...', found $count\"\n  }\n\n  $find2 = $d.Content.Find\n  $find2.Text = $search\n  $find2.Replacement.Text = $replace\n  $find2.Execute($find2.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find2.Replacement.Text, 2) | Out-Null\n}\n\nWrite-Output \"done\"\n"}
